# Update countries & provincias Spain
# Daily refresh of the "Pais" COVID sheet: new timestamp, updated
# case counters for the rows whose totals changed, and a few countries
# that swapped rank (so the country label in that row changes too,
# since the table is kept sorted by total cases descending).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 20 de Julio de 2020 a las 01:57"

# Row 4: Estados Unidos (rank unchanged)
$ws.Range("B4").Value = 3893963
$ws.Range("C4").Value = 60692
$ws.Range("D4").Value = 1802123
$ws.Range("E4").Value = 1948577
$ws.Range("G4").Value = 386
$ws.Range("H4").Value = 143263

# Row 5: Brasil (rank unchanged)
$ws.Range("B5").Value = 2099896
$ws.Range("C5").Value = 24650
$ws.Range("E5").Value = 649134
$ws.Range("G5").Value = 716
$ws.Range("H5").Value = 79533

# Rows 58-59: Irlanda / Kirguistan swap rank
$ws.Range("A58").Value = "Kirguistan"
$ws.Range("B58").Value = 26532
$ws.Range("C58").Value = 1926
$ws.Range("D58").Value = 12328
$ws.Range("E58").Value = 13201
$ws.Range("G58").Value = 103
$ws.Range("H58").Value = 1003

$ws.Range("A59").Value = "Irlanda"
$ws.Range("B59").Value = 25760
$ws.Range("C59").Value = 10
$ws.Range("D59").Value = 23364
$ws.Range("E59").Value = 643
$ws.Range("H59").Value = 1753

# Row 60: Japon (rank unchanged)
$ws.Range("B60").Value = 24642
$ws.Range("C60").Value = 510
$ws.Range("D60").Value = 19576
$ws.Range("E60").Value = 4081

# Rows 74-76: El Salvador / Australia / Venezuela rotate rank
$ws.Range("A74").Value = "Venezuela"
$ws.Range("B74").Value = 11891
$ws.Range("C74").Value = 408
$ws.Range("D74").Value = 3972
$ws.Range("E74").Value = 7807
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 112

$ws.Range("A75").Value = "El Salvador"
$ws.Range("B75").Value = 11846
$ws.Range("C75").Value = 338
$ws.Range("D75").Value = 6705
$ws.Range("E75").Value = 4806
$ws.Range("G75").Value = 11
$ws.Range("H75").Value = 335

$ws.Range("A76").Value = "Australia"
$ws.Range("B76").Value = 11802
$ws.Range("C76").Value = 361
$ws.Range("D76").Value = 8273
$ws.Range("E76").Value = 3407
$ws.Range("G76").Value = 4
$ws.Range("H76").Value = 122

# Rows 143-144: Burkina Faso / Uruguay swap rank
$ws.Range("A143").Value = "Uruguay"
$ws.Range("B143").Value = 1054
$ws.Range("C143").Value = 10
$ws.Range("D143").Value = 922
$ws.Range("E143").Value = 99
$ws.Range("H143").Value = 33

$ws.Range("A144").Value = "Burkina Faso"
$ws.Range("B144").Value = 1052
$ws.Range("C144").Value = 5
$ws.Range("D144").Value = 901
$ws.Range("E144").Value = 98
$ws.Range("H144").Value = 53

# Rows 146-147: Georgia / Surinam swap rank
$ws.Range("A146").Value = "Surinam"
$ws.Range("B146").Value = 1029
$ws.Range("C146").Value = 28
$ws.Range("D146").Value = 627
$ws.Range("E146").Value = 381
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 21

$ws.Range("A147").Value = "Georgia"
$ws.Range("B147").Value = 1028
$ws.Range("C147").Value = 10
$ws.Range("D147").Value = 899
$ws.Range("E147").Value = 114
$ws.Range("H147").Value = 15

# Row 150: Togo (rank unchanged)
$ws.Range("B150").Value = 778
$ws.Range("C150").Value = 4
$ws.Range("D150").Value = 551
$ws.Range("E150").Value = 212

# Row 164: Mauricio (rank unchanged)
$ws.Range("D164").Value = 331
$ws.Range("E164").Value = 2

# Rows 166-168: Isla de Man / Comoras / Guyana rotate rank
$ws.Range("A166").Value = "Guyana"
$ws.Range("C166").Value = 9
$ws.Range("D166").Value = 163
$ws.Range("E166").Value = 154
$ws.Range("H166").Value = 19

$ws.Range("A167").Value = "Isla de Man"
$ws.Range("B167").Value = 336
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 312
$ws.Range("E167").Value = 0
$ws.Range("H167").Value = 24

$ws.Range("A168").Value = "Comoras"
$ws.Range("B168").Value = 334
$ws.Range("C168").Value = 6
$ws.Range("D168").Value = 313
$ws.Range("E168").Value = 14
$ws.Range("H168").Value = 7

# Rows 210-211: Islas Malvinas / Groenlandia swap rank (counts identical)
$ws.Range("A210").Value = "Groenlandia"

$ws.Range("A211").Value = "Islas Malvinas"
